# "Added Full BOM and gerber"
# The BOM's Designator column (C) originally held the bus/signal names
# feeding each connector; update it to the actual reference designators
# (J1-J8) used on the gerber/PCB silkscreen.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")
$ws.Activate()

$ws.Range("C2").Value  = "J1,J2,J3,J4"   # was I2C2,USART2,I2C1,USART1
$ws.Range("C9").Value  = "J6"            # was SPI
$ws.Range("C22").Value = "J5"            # was JTAG
$ws.Range("C23").Value = "J7"            # was USB
$ws.Range("C26").Value = "J8"            # was VEXT

$ws.Range("C26").Select()
